$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = @'
questions = [
    {
        "title": "You are developing a Next.js website on a laptop and want to preview how it will look on mobile. To do this, you need to make your Next.js app accessible via the local area network IP address 192.168.1.2. This will allow you to access the development version of the website directly on your phone.Which Next.js CLI command should you use to achieve this?",
        "ques_type": 2,
        "options": [
            "npx next dev --hostname 192.168.1.2",
            "npx next dev -hostname 192.168.1.2",
            "npx next dev -h 192.168.1.2",
            "npx next dev --H 192.168.1.2"
        ],
        "score": "npx next dev --hostname 192.168.1.2"
    },
    {
        "title": "You work for a startup trying to reach more users through blog articles that introduce its products to readers. You're developing the blog using Next.js, with content coming from a headless content management system (CMS). Your project manager wants the blog to be easily indexed by search engines and quickly accessed by users so there are no long loading times when moving from one article to another.Which data-fetching method should you use?",
        "ques_type": 2,
        "options": [
            "getServerSideProps",
            "getStaticPaths",
            "getStaticProps",
            "getInitialProps"
        ],
        "score": "getStaticProps"
    },
    {
        "title": "You work for an ecommerce company that uses Next.js. You have completed the product list page and are developing a product detail page to display the details of each product by product ID.Which of the following file names should you use to create the page?",
        "ques_type": 2,
        "options": [
            "product.js",
            "[product-id].js",
            "product-id.js",
            "id.js"
        ],
        "score": "[product-id].js"
    },
    {
        "title": "Your website has a component called MobileNav, which appears when mobile users scroll. To improve initial loading performance, you plan to use code splitting with dynamic import. Which of the following codes should you use to import MobileNav dynamically?",
        "ques_type": 2,
        "options": [
            "const MobileNav = dynamic(() =&gt import('../components/MobileNav'))",
            "const MobileNav = import(() =&gt dynamic('../components/MobileNav'))",
            "const MobileNav = dynamicImport(() =&gt dynamic('../components/MobileNav'))",
            "const MobileNav = import(() =&gt dynamicImport('../components/MobileNav'))"
        ],
        "score": "const MobileNav = dynamic(() =&gt import('../components/MobileNav'))"
    }
]
'@

$ws.Range("A2").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
$ws.Range("A1").Style = "Normal"
$ws.Range("A1").Value = $text
$ws.Rows.Item(1).AutoFit()
